$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Spring")
$ws2 = $wb.Worksheets.Item("Summer")

# ---------------------------------------------------------------------------
# Create the four brand-new shared strings in the exact order they first
# appear in the target workbook (RANGE, INTERVALS, "total q ", "just samples")
# so they land at shared-string indices 104-107 respectively.
# ---------------------------------------------------------------------------
$ws1.Range("F17").Value = "RANGE"
$ws1.Range("F18").Value = "INTERVALS"
$ws1.Range("F11").Value = "total q "
$ws1.Range("F1").Value = "just samples"

# ---------------------------------------------------------------------------
# Sheet "Spring" - little "range / intervals" helper table in columns F/G.
# Enter the G17/G18 formulas first (while their precedents are still blank)
# so they don't pick up the 0.0000 number format that G12-G15 get below.
# ---------------------------------------------------------------------------
$ws1.Range("G17").Formula = "=G15-G12"
$ws1.Range("G18").Formula = "=G17/3"

$ws1.Range("F12").Value = "MIN"
$ws1.Range("G12").Value = 0.012
$ws1.Range("G12").NumberFormat = "0.0000"

$ws1.Range("F13").Value = "T1"
$ws1.Range("G13").Formula = "=G12+G18"
$ws1.Range("G13").NumberFormat = "0.0000"

$ws1.Range("F14").Value = "T2"
$ws1.Range("G14").Formula = "=G13+G18"
$ws1.Range("G14").NumberFormat = "0.0000"

$ws1.Range("F15").Value = "MAX"
$ws1.Range("G15").Value = 0.298
$ws1.Range("G15").NumberFormat = "0.0000"

# ---------------------------------------------------------------------------
# Sheet "Summer" - same helper table, shifted up one row compared to Spring.
# ---------------------------------------------------------------------------
$ws2.Range("F1").Value = "just samples"
$ws2.Range("F1").HorizontalAlignment = -4108

$ws2.Range("F10").Value = "total q "

$ws2.Range("F16").Value = "RANGE"
$ws2.Range("G16").Formula = "=G14-G11"

$ws2.Range("F17").Value = "INTERVALS"
$ws2.Range("G17").Formula = "=G16/3"

$ws2.Range("F11").Value = "MIN"
$ws2.Range("G11").Value = 0.002

$ws2.Range("F12").Value = "T1"
$ws2.Range("G12").Formula = "=G11+G17"

$ws2.Range("F13").Value = "T2"
$ws2.Range("G13").Formula = "=G12+G17"

$ws2.Range("F14").Value = "MAX"
$ws2.Range("G14").Value = 0.098

# ---------------------------------------------------------------------------
# Restore the selections recorded in each sheet view.
# ---------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("I12").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("I16").Select() | Out-Null
